$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.887.68"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "1.617.30"
$ws.Range("E4").Value = "  -0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.39"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.26"
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "1.841.94"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").Value = "1.616.48"
$ws.Range("E13").Value = "  -4.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.12"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "25.889.03"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.40"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.02"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.48"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.01"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.62"
$ws.Range("E28").Value = "  -2.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.22"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("D36").Value = "1.129.07"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  -4.62%  "
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.510"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.88"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "1.753.59"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.748"
$ws.Range("E43").Value = "  -4.85%  "
$ws.Range("E44").Value = "  -4.29%  "
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.91"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.411"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.45"
$ws.Range("E51").Value = "  -1.86%  "
